$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) <meta name="title" ...> : "1 What is Blender" -> " 4 Making Curves with Loop Tools"
# ------------------------------------------------------------------
$p = $d.Paragraphs.Item(1)
$p.Range.Find.Execute("1 What is Blender", $true, $false, $false, $false, $false, $true, 1, $false, " 4 Making Curves with Loop Tools", 2) | Out-Null

# ------------------------------------------------------------------
# 2) <meta name="keywords" ...> : prepend "Making Curves with Loop Tools, "
#    in front of the existing "Blender, 3D Modeling, Animation, Graphic Art"
# ------------------------------------------------------------------
$p = $d.Paragraphs.Item(6)
$p.Range.Find.Execute("Blender, 3D Modeling, Animation, Graphic Art", $true, $false, $false, $false, $false, $true, 1, $false, "Making Curves with Loop Tools, Blender, 3D Modeling, Animation, Graphic Art", 2) | Out-Null

# ------------------------------------------------------------------
# 3) <meta name="description" ...> : rewrite what the article explains
# ------------------------------------------------------------------
$p = $d.Paragraphs.Item(9)
$oldDescription = "what the 3D modeling program " + [char]0x201C + "Blender " + [char]0x22 + " is all about./>"
$newDescription = "how to go about making Curves with Loop Tools />"
$p.Range.Find.Execute($oldDescription, $true, $false, $false, $false, $false, $true, 1, $false, $newDescription, 2) | Out-Null

# ------------------------------------------------------------------
# 4) <meta name="category" ...> : prepend "Making Curves with Loop Tools, "
#    (wrapped in a bookmark, as Word does for pasted/auto-linked text)
#    in front of the existing "Blender, 3D Modeling, Animation, Graphic Art"
# ------------------------------------------------------------------
$p = $d.Paragraphs.Item(11)
$p.Range.Find.Execute("Blender, 3D Modeling, Animation, Graphic Art", $true, $false, $false, $false, $false, $true, 1, $false, "Making Curves with Loop Tools, Blender, 3D Modeling, Animation, Graphic Art", 2) | Out-Null

$p = $d.Paragraphs.Item(11)
$rng = $p.Range
$rng.Find.Execute("Making Curves with Loop Tools", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$d.Bookmarks.Add("_Hlk190486324", $rng) | Out-Null

# ------------------------------------------------------------------
# 5) <meta name="revised" ...> : bump the revision date
# ------------------------------------------------------------------
$p = $d.Paragraphs.Item(18)
$p.Range.Find.Execute("Wednesday, December 11, 2024", $true, $false, $false, $false, $false, $true, 1, $false, "Saturday, February 15, 2025", 2) | Out-Null

# ------------------------------------------------------------------
# 6) <meta name="url" ...> : point at the new article location
# ------------------------------------------------------------------
$p = $d.Paragraphs.Item(20)
$oldUrl = "Enlightenment/Articles/2024/8-Blender-2024/1-What-Is-Blender/1-What-Is-Blender.html"
$newUrl = "Enlightenment/Articles/2025/1-Blender-Continued/7-Loop-Tools/4-Curves/4-Curves.html"
$p.Range.Find.Execute($oldUrl, $true, $false, $false, $false, $false, $true, 1, $false, $newUrl, 2) | Out-Null
